$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..G, rows 2..9 (filtered sval data regen)
$data = @{
    2 = @(0.6606524410359556, 3286.919754855326, 0.1494219747398047, 10.19245300693656, 1, 3297.922282278038)
    3 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 0, 12.65830550497259)
    4 = @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 0, 12.1040582534506)
    5 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 2.960089034096801)
    6 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 3.754798637575387)
    7 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    8 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    9 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
